{"js": "// Apply the CSS exercises 1-7 copy-edit pass to the HTML notes document.\n// The real-world commit mostly re-flowed existing runs (proofing marks,\n// run splits) with no visible text impact; the only reader-visible\n// changes are five small wording/formatting fixes inside the \"img/video/\n// audio/icon src path\" placeholders and the <img alt> sentence. Apply\n// each one with a scoped search + replace so we don't disturb anything\n// else in the document.\nconst body = context.document.body;\n\nconst replacements = [\n  // \"<img src=\u201dRuta de la imagen\u201d ...>\" -> spaces in the placeholder become underscores\n  [\" de la imagen\u201d\", \"_de_la_imagen\u201d\"],\n  // \"...alt=\u201dTexto en caso que no se cargue\u201d>\" -> \"caso que\" becomes \"caso de que\"\n  [\"caso que no se cargue\", \"caso de que no se cargue\"],\n  // \"<video src=\u201dRuta del v\u00eddeo\u201d ...>\" -> spaces become underscores\n  [\"Ruta del v\u00eddeo\u201d\", \"Ruta_del_v\u00eddeo\u201d\"],\n  // \"<audio src=\u201dRuta del audio\u201d ...>\" -> spaces become underscores\n  [\"Ruta del audio\u201d\", \"Ruta_del_audio\u201d\"],\n  // \"<link rel=\u201dicon\u201d href=\u201dRuta del \u00edcono\u201d>\" -> spaces become underscores\n  [\"Ruta del \u00edcono\u201d\", \"Ruta_del_\u00edcono\u201d\"],\n];\n\nfor (const [search, replacement] of replacements) {\n  const found = body.search(search, { matchCase: true, matchWholeWord: false });\n  found.load(\"text\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    continue;\n  }\n\n  found.items[0].insertText(replacement, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Apply the CSS exercises 1-7 copy-edit pass to the HTML notes document.\n# The real-world commit mostly re-flowed existing runs (proofing marks,\n# run splits) with no visible text impact; the only reader-visible\n# changes are five small wording/formatting fixes inside the \"img/video/\n# audio/icon src path\" placeholders and the <img alt> sentence. Apply\n# each one with a scoped Find/Replace so we don't disturb anything else\n# in the document.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\" de la imagen\u201d\", \"_de_la_imagen\u201d\"),\n    @(\"caso que no se cargue\", \"caso de que no se cargue\"),\n    @(\"Ruta del v\u00eddeo\u201d\", \"Ruta_del_v\u00eddeo\u201d\"),\n    @(\"Ruta del audio\u201d\", \"Ruta_del_audio\u201d\"),\n    @(\"Ruta del \u00edcono\u201d\", \"Ruta_del_\u00edcono\u201d\")\n)\n\nforeach ($pair in $replacements) {\n    $searchText = $pair[0]\n    $replaceText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $searchText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 1\n\n    $find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
